$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 60561.824
$ws.Range("I98").Value = 84829.25
$ws.Range("J98").Value = 2320
$ws.Range("K98").Value = 84829.25
$ws.Range("L98").Value = 2320
$ws.Range("M98").Value = -83331.25
$ws.Range("N98").Value = -5316

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 770
$ws.Range("I115").Value = 770
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2310
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -743
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 60561.824
$ws.Range("I122").Value = 84829.25
$ws.Range("J122").Value = 2320
$ws.Range("K122").Value = 254487.75
$ws.Range("L122").Value = 6960
$ws.Range("M122").Value = -252037.75
$ws.Range("N122").Value = -11860

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1800.2778
$ws.Range("I132").Value = 1212.8125
$ws.Range("K132").Value = 3638.4375
$ws.Range("M132").Value = -1108.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1397.421
$ws.Range("I45").Value = 1287.8
$ws.Range("J45").Value = 1808.5
$ws.Range("K45").Value = 1287.8
$ws.Range("L45").Value = 1808.5
$ws.Range("M45").Value = -910.8
$ws.Range("N45").Value = -2562.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2018.2667
$ws.Range("I61").Value = 1632.4
$ws.Range("J61").Value = 2790
$ws.Range("K61").Value = 1632.4
$ws.Range("L61").Value = 2790
$ws.Range("M61").Value = -1420.4
$ws.Range("N61").Value = -3214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1379.4667
$ws.Range("I74").Value = 922.2222
$ws.Range("J74").Value = 2065.3333
$ws.Range("K74").Value = 922.2222
$ws.Range("L74").Value = 2065.3333
$ws.Range("M74").Value = -48.22220000000004
$ws.Range("N74").Value = -3813.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1379.4667
$ws.Range("I77").Value = 922.2222
$ws.Range("J77").Value = 2065.3333
$ws.Range("K77").Value = 4611.111
$ws.Range("L77").Value = 10326.6665
$ws.Range("M77").Value = -243.1109999999999
$ws.Range("N77").Value = -19062.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2018.2667
$ws.Range("I136").Value = 1632.4
$ws.Range("J136").Value = 2790
$ws.Range("K136").Value = 4897.200000000001
$ws.Range("L136").Value = 8370
$ws.Range("M136").Value = -2347.200000000001
$ws.Range("N136").Value = -13470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 35641.25
$ws.Range("I88").Value = 30000
$ws.Range("J88").Value = 41282.5
$ws.Range("K88").Value = 30000
$ws.Range("L88").Value = 41282.5
$ws.Range("M88").Value = -29594
$ws.Range("N88").Value = -42094.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 35641.25
$ws.Range("I91").Value = 30000
$ws.Range("J91").Value = 41282.5
$ws.Range("K91").Value = 30000
$ws.Range("L91").Value = 41282.5
$ws.Range("M91").Value = -28596
$ws.Range("N91").Value = -44090.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1685.2106
$ws.Range("I99").Value = 1843.25
$ws.Range("J99").Value = 1414.2858
$ws.Range("K99").Value = 1843.25
$ws.Range("L99").Value = 1414.2858
$ws.Range("M99").Value = -345.25
$ws.Range("N99").Value = -4410.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1920.7667
$ws.Range("I31").Value = 1717.2307
$ws.Range("J31").Value = 3243.75
$ws.Range("K31").Value = 1717.2307
$ws.Range("L31").Value = 3243.75
$ws.Range("M31").Value = -1422.2307
$ws.Range("N31").Value = -3833.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1920.7667
$ws.Range("I34").Value = 1717.2307
$ws.Range("J34").Value = 3243.75
$ws.Range("K34").Value = 1717.2307
$ws.Range("L34").Value = 3243.75
$ws.Range("M34").Value = -1515.2307
$ws.Range("N34").Value = -3647.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 39666.668
$ws.Range("J88").Value = 39666.668
$ws.Range("L88").Value = 39666.668
$ws.Range("N88").Value = -40478.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 39666.668
$ws.Range("J91").Value = 39666.668
$ws.Range("L91").Value = 39666.668
$ws.Range("N91").Value = -42474.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3142.8572
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 940.2917
$ws.Range("I107").Value = 585.2
$ws.Range("J107").Value = 1532.1111
$ws.Range("K107").Value = 585.2
$ws.Range("L107").Value = 1532.1111
$ws.Range("M107").Value = 1334.8
$ws.Range("N107").Value = -5372.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1376.4286
$ws.Range("I122").Value = 1454
$ws.Range("J122").Value = 911
$ws.Range("K122").Value = 4362
$ws.Range("L122").Value = 2733
$ws.Range("M122").Value = -1912
$ws.Range("N122").Value = -7633

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3142.8572
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2179.3076
$ws.Range("I132").Value = 677.375
$ws.Range("J132").Value = 4582.4
$ws.Range("K132").Value = 2032.125
$ws.Range("L132").Value = 13747.2
$ws.Range("M132").Value = 497.875
$ws.Range("N132").Value = -18807.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1732.6875
$ws.Range("I134").Value = 1379.2142
$ws.Range("J134").Value = 4207
$ws.Range("K134").Value = 4137.642599999999
$ws.Range("L134").Value = 12621
$ws.Range("M134").Value = -1602.642599999999
$ws.Range("N134").Value = -17691

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5366.6665
$ws.Range("I133").Value = 4866.6665
$ws.Range("K133").Value = 14599.9995
$ws.Range("M133").Value = -9539.999500000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 503.04544
$ws.Range("I22").Value = 606.6667
$ws.Range("J22").Value = 378.7
$ws.Range("K22").Value = 606.6667
$ws.Range("L22").Value = 378.7
$ws.Range("M22").Value = -311.6667
$ws.Range("N22").Value = -968.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 503.04544
$ws.Range("I27").Value = 606.6667
$ws.Range("J27").Value = 378.7
$ws.Range("K27").Value = 606.6667
$ws.Range("L27").Value = 378.7
$ws.Range("M27").Value = -499.6667
$ws.Range("N27").Value = -592.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 23271
$ws.Range("J69").Value = 23271
$ws.Range("L69").Value = 23271
$ws.Range("N69").Value = -24769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 23271
$ws.Range("J72").Value = 23271
$ws.Range("L72").Value = 69813
$ws.Range("N72").Value = -77301
